# Update "想去人数" (F column) figures across the four sheets to reflect
# refreshed scrape counts, as published to gh-pages.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 141
$ws.Range("F3").Value = 1326
$ws.Range("F4").Value = 1124
$ws.Range("F5").Value = 1015
$ws.Range("F6").Value = 1792
$ws.Range("F7").Value = 557
$ws.Range("F8").Value = 1198
$ws.Range("F12").Value = 293
$ws.Range("F13").Value = 66
$ws.Range("F16").Value = 167
$ws.Range("F20").Value = 328
$ws.Range("F21").Value = 144
$ws.Range("F22").Value = 669
$ws.Range("F23").Value = 33
$ws.Range("F24").Value = 644
$ws.Range("F27").Value = 874
$ws.Range("F28").Value = 312
$ws.Range("F29").Value = 160
$ws.Range("F30").Value = 43

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 318
$ws.Range("F7").Value = 254
$ws.Range("F9").Value = 3
$ws.Range("F11").Value = 120

# ---- Sheet "本地生活" ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 309

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 309
$ws.Range("F3").Value = 141
$ws.Range("F4").Value = 1326
$ws.Range("F5").Value = 1124
$ws.Range("F6").Value = 1015
$ws.Range("F7").Value = 1792
$ws.Range("F8").Value = 557
$ws.Range("F9").Value = 1198
$ws.Range("F14").Value = 293
$ws.Range("F15").Value = 66
$ws.Range("F18").Value = 167
$ws.Range("F22").Value = 318
$ws.Range("F25").Value = 328
$ws.Range("F27").Value = 254
$ws.Range("F28").Value = 254
$ws.Range("F29").Value = 144
$ws.Range("F30").Value = 669
$ws.Range("F31").Value = 33
$ws.Range("F32").Value = 644
$ws.Range("F35").Value = 874
$ws.Range("F36").Value = 312
$ws.Range("F38").Value = 3
$ws.Range("F39").Value = 160
$ws.Range("F40").Value = 43
$ws.Range("F43").Value = 120
$ws.Range("F44").Value = 120

$wb.Save()
